# Generate Report for Archive
# Update the Status of the two files that have moved from "Ready for handoff"
# to "In Translation" (811bc1ec-... and 9b9a0341-...) across every sheet:
# Overview (columns B & C) and the per-locale sheets zh-cn / de-de (column C).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---------------------------------------------------
$ws = $wb.Worksheets.Item("Overview")
$ws.Range("B3").Value = "In Translation"
$ws.Range("C3").Value = "In Translation"
$ws.Range("B4").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"

# --- zh-cn sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("zh-cn")
$ws.Range("C3").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"

# --- de-de sheet --------------------------------------------------------
$ws = $wb.Worksheets.Item("de-de")
$ws.Range("C3").Value = "In Translation"
$ws.Range("C4").Value = "In Translation"
